$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Angpt1"
$ws.Cells.Item(2, 3).Value = "Itga5"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.110507
$ws.Cells.Item(2, 8).Value = 0.331521
$ws.Cells.Item(2, 9).Value = 0.004605687348208628
$ws.Cells.Item(2, 10).Value = 0.004605687348208628
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 30.99161333333333
$ws.Cells.Item(2, 14).Value = 92.97484
$ws.Cells.Item(2, 15).Value = 0.3599121977633812
$ws.Cells.Item(2, 16).Value = 0.3599121977633811
$ws.Cells.Item(2, 17).Value = 3.424790214626667
$ws.Cells.Item(2, 18).Value = 30.82311193164
$ws.Cells.Item(2, 19).Value = 0.001657643055704766
$ws.Cells.Item(2, 20).Value = 0.001657643055704766

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Angpt1"
$ws.Cells.Item(3, 3).Value = "Itga5"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.110507
$ws.Cells.Item(3, 8).Value = 0.331521
$ws.Cells.Item(3, 9).Value = 0.004605687348208628
$ws.Cells.Item(3, 10).Value = 0.004605687348208628
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 29.913269
$ws.Cells.Item(3, 14).Value = 89.739807
$ws.Cells.Item(3, 15).Value = 0.3473891556493311
$ws.Cells.Item(3, 16).Value = 0.3473891556493311
$ws.Cells.Item(3, 17).Value = 3.305625617383
$ws.Cells.Item(3, 18).Value = 29.750630556447
$ws.Cells.Item(3, 19).Value = 0.001599965839079002
$ws.Cells.Item(3, 20).Value = 0.001599965839079002

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Angpt1"
$ws.Cells.Item(4, 3).Value = "Itga5"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.110507
$ws.Cells.Item(4, 8).Value = 0.331521
$ws.Cells.Item(4, 9).Value = 0.004605687348208628
$ws.Cells.Item(4, 10).Value = 0.004605687348208628
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 19.150218
$ws.Cells.Item(4, 14).Value = 57.450654
$ws.Cells.Item(4, 15).Value = 0.2223955550134164
$ws.Cells.Item(4, 16).Value = 0.2223955550134163
$ws.Cells.Item(4, 17).Value = 2.116233140526
$ws.Cells.Item(4, 18).Value = 19.046098264734
$ws.Cells.Item(4, 19).Value = 0.001024284394023128
$ws.Cells.Item(4, 20).Value = 0.001024284394023127

$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Angpt1"
$ws.Cells.Item(5, 3).Value = "Itga5"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.110507
$ws.Cells.Item(5, 8).Value = 0.331521
$ws.Cells.Item(5, 9).Value = 0.004605687348208628
$ws.Cells.Item(5, 10).Value = 0.004605687348208628
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 6.053716000000001
$ws.Cells.Item(5, 14).Value = 18.161148
$ws.Cells.Item(5, 15).Value = 0.07030309157387134
$ws.Cells.Item(5, 16).Value = 0.07030309157387132
$ws.Cells.Item(5, 17).Value = 0.6689779940120001
$ws.Cells.Item(5, 18).Value = 6.020801946108
$ws.Cells.Item(5, 19).Value = 0.0003237940594017318
$ws.Cells.Item(5, 20).Value = 0.0003237940594017317

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Angpt1"
$ws.Cells.Item(6, 3).Value = "Itga5"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 22.14783133333333
$ws.Cells.Item(6, 8).Value = 66.443494
$ws.Cells.Item(6, 9).Value = 0.9230726249214253
$ws.Cells.Item(6, 10).Value = 0.9230726249214253
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 30.99161333333333
$ws.Cells.Item(6, 14).Value = 92.97484
$ws.Cells.Item(6, 15).Value = 0.3599121977633812
$ws.Cells.Item(6, 16).Value = 0.3599121977633811
$ws.Cells.Item(6, 17).Value = 686.3970248545511
$ws.Cells.Item(6, 18).Value = 6177.57322369096
$ws.Cells.Item(6, 19).Value = 0.3322250971306834
$ws.Cells.Item(6, 20).Value = 0.3322250971306833

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Angpt1"
$ws.Cells.Item(7, 3).Value = "Itga5"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 22.14783133333333
$ws.Cells.Item(7, 8).Value = 66.443494
$ws.Cells.Item(7, 9).Value = 0.9230726249214253
$ws.Cells.Item(7, 10).Value = 0.9230726249214253
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 29.913269
$ws.Cells.Item(7, 14).Value = 89.739807
$ws.Cells.Item(7, 15).Value = 0.3473891556493311
$ws.Cells.Item(7, 16).Value = 0.3473891556493311
$ws.Cells.Item(7, 17).Value = 662.5140364406286
$ws.Cells.Item(7, 18).Value = 5962.626327965658
$ws.Cells.Item(7, 19).Value = 0.3206654197744657
$ws.Cells.Item(7, 20).Value = 0.3206654197744656

$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Angpt1"
$ws.Cells.Item(8, 3).Value = "Itga5"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 22.14783133333333
$ws.Cells.Item(8, 8).Value = 66.443494
$ws.Cells.Item(8, 9).Value = 0.9230726249214253
$ws.Cells.Item(8, 10).Value = 0.9230726249214253
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 19.150218
$ws.Cells.Item(8, 14).Value = 57.450654
$ws.Cells.Item(8, 15).Value = 0.2223955550134164
$ws.Cells.Item(8, 16).Value = 0.2223955550134163
$ws.Cells.Item(8, 17).Value = 424.1357982605639
$ws.Cells.Item(8, 18).Value = 3817.222184345076
$ws.Cells.Item(8, 19).Value = 0.2052872487370915
$ws.Cells.Item(8, 20).Value = 0.2052872487370915

$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Angpt1"
$ws.Cells.Item(9, 3).Value = "Itga5"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 22.14783133333333
$ws.Cells.Item(9, 8).Value = 66.443494
$ws.Cells.Item(9, 9).Value = 0.9230726249214253
$ws.Cells.Item(9, 10).Value = 0.9230726249214253
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 6.053716000000001
$ws.Cells.Item(9, 14).Value = 18.161148
$ws.Cells.Item(9, 15).Value = 0.07030309157387134
$ws.Cells.Item(9, 16).Value = 0.07030309157387132
$ws.Cells.Item(9, 17).Value = 134.0766809079013
$ws.Cells.Item(9, 18).Value = 1206.690128171112
$ws.Cells.Item(9, 19).Value = 0.06489485927918476
$ws.Cells.Item(9, 20).Value = 0.06489485927918474

$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Angpt1"
$ws.Cells.Item(10, 3).Value = "Itga5"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.735257333333333
$ws.Cells.Item(10, 8).Value = 5.205772
$ws.Cells.Item(10, 9).Value = 0.07232168773036617
$ws.Cells.Item(10, 10).Value = 0.07232168773036617
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 30.99161333333333
$ws.Cells.Item(10, 14).Value = 92.97484
$ws.Cells.Item(10, 15).Value = 0.3599121977633812
$ws.Cells.Item(10, 16).Value = 0.3599121977633811
$ws.Cells.Item(10, 17).Value = 53.77842430849778
$ws.Cells.Item(10, 18).Value = 484.0058187764799
$ws.Cells.Item(10, 19).Value = 0.02602945757699305
$ws.Cells.Item(10, 20).Value = 0.02602945757699304

$ws.Cells.Item(11, 1).Value = "sCs"
$ws.Cells.Item(11, 2).Value = "Angpt1"
$ws.Cells.Item(11, 3).Value = "Itga5"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 1.735257333333333
$ws.Cells.Item(11, 8).Value = 5.205772
$ws.Cells.Item(11, 9).Value = 0.07232168773036617
$ws.Cells.Item(11, 10).Value = 0.07232168773036617
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 29.913269
$ws.Cells.Item(11, 14).Value = 89.739807
$ws.Cells.Item(11, 15).Value = 0.3473891556493311
$ws.Cells.Item(11, 16).Value = 0.3473891556493311
$ws.Cells.Item(11, 17).Value = 51.90721939622266
$ws.Cells.Item(11, 18).Value = 467.1649745660039
$ws.Cells.Item(11, 19).Value = 0.0251237700357865
$ws.Cells.Item(11, 20).Value = 0.02512377003578649

$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Angpt1"
$ws.Cells.Item(12, 3).Value = "Itga5"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 1.735257333333333
$ws.Cells.Item(12, 8).Value = 5.205772
$ws.Cells.Item(12, 9).Value = 0.07232168773036617
$ws.Cells.Item(12, 10).Value = 0.07232168773036617
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 19.150218
$ws.Cells.Item(12, 14).Value = 57.450654
$ws.Cells.Item(12, 15).Value = 0.2223955550134164
$ws.Cells.Item(12, 16).Value = 0.2223955550134163
$ws.Cells.Item(12, 17).Value = 33.230556219432
$ws.Cells.Item(12, 18).Value = 299.075005974888
$ws.Cells.Item(12, 19).Value = 0.01608402188230177
$ws.Cells.Item(12, 20).Value = 0.01608402188230177

$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Angpt1"
$ws.Cells.Item(13, 3).Value = "Itga5"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 1.735257333333333
$ws.Cells.Item(13, 8).Value = 5.205772
$ws.Cells.Item(13, 9).Value = 0.07232168773036617
$ws.Cells.Item(13, 10).Value = 0.07232168773036617
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 6.053716000000001
$ws.Cells.Item(13, 14).Value = 18.161148
$ws.Cells.Item(13, 15).Value = 0.07030309157387134
$ws.Cells.Item(13, 16).Value = 0.07030309157387132
$ws.Cells.Item(13, 17).Value = 10.50475508291733
$ws.Cells.Item(13, 18).Value = 94.54279574625599
$ws.Cells.Item(13, 19).Value = 0.005084438235284859
$ws.Cells.Item(13, 20).Value = 0.005084438235284859
